$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09102500000000001
$ws.Range("H2").Value = 0.273075
$ws.Range("I2").Value = 0.8515285885346505
$ws.Range("J2").Value = 0.8515285885346504
$ws.Range("M2").Value = 0.08268033333333334
$ws.Range("Q2").Value = 0.007525977341666669
$ws.Range("R2").Value = 0.06773379607500001
$ws.Range("S2").Value = 0.8515285885346505
$ws.Range("T2").Value = 0.8515285885346504

# Row 3 updates
$ws.Range("I3").Value = 0.1484714114653495
$ws.Range("J3").Value = 0.1484714114653495
$ws.Range("M3").Value = 0.08268033333333334
$ws.Range("S3").Value = 0.1484714114653495
$ws.Range("T3").Value = 0.1484714114653495

$wb.Save()
